$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column D (experimentDesign) to bound the replace
$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row

# Replace "Environmental Perturbation" -> "Environmental_Perturbation" in column D
$rngD = $ws.Range("D1:D$lastRow")
$rngD.Replace("Environmental Perturbation", "Environmental_Perturbation", 1, 1, $false, $false, $false)

# Replace "KN99 alpha" -> "KN99_alpha" in column F
$rngF = $ws.Range("F1:F$lastRow")
$rngF.Replace("KN99 alpha", "KN99_alpha", 1, 1, $false, $false, $false)

# Update the selection to F2:F27 with active cell F2
$ws.Range("F2:F27").Select()
